# Reorder the worker/period rows (B16:J27) in Hoja1 grouping by employee:
# first all JOSE DAVID GUERRA MEJIA periods (2504,2503,2502,2501,2412,2411),
# then all DIEGO ANDRES ARRIETA BOHORQUEZ periods (2504,2503,2502,2501,2412,2411).
# Rows 28-29 are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$docId_Jose = "92191907"
$name_Jose  = "JOSE DAVID GUERRA MEJIA"
$docId_Diego = "1235046712"
$name_Diego  = "DIEGO ANDRES ARRIETA BOHORQUEZ"

$rows = @(
    @{ Row = 16; DocId = $docId_Jose;  Name = $name_Jose;  Periodo = "2504"; Mora = 60000; Salario = 1500000 },
    @{ Row = 17; DocId = $docId_Jose;  Name = $name_Jose;  Periodo = "2503"; Mora = 60000; Salario = 1500000 },
    @{ Row = 18; DocId = $docId_Jose;  Name = $name_Jose;  Periodo = "2502"; Mora = 60000; Salario = 1500000 },
    @{ Row = 19; DocId = $docId_Jose;  Name = $name_Jose;  Periodo = "2501"; Mora = 60000; Salario = 1500000 },
    @{ Row = 20; DocId = $docId_Jose;  Name = $name_Jose;  Periodo = "2412"; Mora = 60000; Salario = 1500000 },
    @{ Row = 21; DocId = $docId_Jose;  Name = $name_Jose;  Periodo = "2411"; Mora = 60000; Salario = 1500000 },
    @{ Row = 22; DocId = $docId_Diego; Name = $name_Diego; Periodo = "2504"; Mora = 52000; Salario = 1300000 },
    @{ Row = 23; DocId = $docId_Diego; Name = $name_Diego; Periodo = "2503"; Mora = 52000; Salario = 1300000 },
    @{ Row = 24; DocId = $docId_Diego; Name = $name_Diego; Periodo = "2502"; Mora = 52000; Salario = 1300000 },
    @{ Row = 25; DocId = $docId_Diego; Name = $name_Diego; Periodo = "2501"; Mora = 52000; Salario = 1300000 },
    @{ Row = 26; DocId = $docId_Diego; Name = $name_Diego; Periodo = "2412"; Mora = 52000; Salario = 1300000 },
    @{ Row = 27; DocId = $docId_Diego; Name = $name_Diego; Periodo = "2411"; Mora = 52000; Salario = 1300000 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 3).Value = $r.DocId     # C: N Doc Trabajador
    $ws.Cells.Item($row, 4).Value = $r.Name      # D: Nombre Trabajador
    $ws.Cells.Item($row, 5).Value = $r.Periodo   # E: Periodo Mora
    $ws.Cells.Item($row, 6).Value = $r.Mora      # F: Valor Mora
    $ws.Cells.Item($row, 7).Value = $r.Salario   # G: Salario Basico
}
